# Add results for 10^5 and 10^6 DT table size
#
# The workbook has two sheets:
#   ~3e5  (Worksheets.Item(1)) - results for a 10^5-row DT table
#   ~3e6  (Worksheets.Item(2)) - results for a 10^6-row DT table
#
# This script fills in the previously-missing "Partial denormalization"
# column (F) in the ~3e6 sheet's second results table, corrects the first
# table's F13 figure, sets that sheet's print setup, and leaves the
# workbook with the ~3e5 sheet on top (selection/active-tab housekeeping
# that comes along with the edit).

$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item(1)   # ~3e5
$ws6 = $wb.Worksheets.Item(2)   # ~3e6

# --- ~3e6 sheet: update the "Partial denormalization" time figure ---
$ws6.Range("F13").Value = 0.61

# --- ~3e6 sheet: fill in the newly-measured "Partial denormalization"
#     numbers for the second table (rows 25-28); G column recalculates
#     automatically from the existing formulas ---
$ws6.Range("F25").Value = 1.36
$ws6.Range("F26").Value = 5.41
$ws6.Range("F27").Value = 3.71
$ws6.Range("F28").Value = 2.35

# --- ~3e6 sheet: page setup for printing ---
$ws6.PageSetup.PaperSize = 9      # xlPaperA4
$ws6.PageSetup.Orientation = 2    # xlLandscape

# --- leave the selection on ~3e6 where work just happened ... ---
[void]$ws6.Range("F8").Select()

# --- ...then switch back to ~3e5 as the active/top sheet with its own
#     selection, matching the saved workbook view state ---
[void]$ws5.Activate()
[void]$ws5.Range("H52").Select()
